$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "dog": append a new test row (row 6) with the same per-column look
# as the row above it, then fill in the new values.
# ---------------------------------------------------------------------------
$dog = $wb.Worksheets.Item("dog")

$dog.Range("A5:L5").Copy()
$dog.Range("A6:L6").PasteSpecial(-4122)  # xlPasteFormats

$dog.Range("A6").Value = 45781
$dog.Range("B6").Value = "PRESENCE"
$dog.Range("C6").Value = 0.47222222222222221
$dog.Range("D6").Value = 0.67708333333333337
$dog.Range("E6").Value = 14
$dog.Range("F6").Value = 4
$dog.Range("G6").Value = "Sunny, cool"
$dog.Range("H6").Value = $true
$dog.Range("I6").Value = "18 minutes 48 seconds"
$dog.Range("J6").Value = 1126
$dog.Range("K6").Value = "Secondary checks"
$dog.Range("L6").Value = "Worked downill. Missed on transects but picked up when returned to point of interest at logs. Koda a bit distracted."

$dog.Range("L16").Select()

# ---------------------------------------------------------------------------
# Sheet "human": add a new leading "Searcher" column, filled in for every
# existing row.
# ---------------------------------------------------------------------------
$human = $wb.Worksheets.Item("human")

$human.Range("A1").EntireColumn.Insert()

$human.Range("A1").Value = "Searcher"
$human.Range("A1").Font.Bold = $true
$human.Range("A1").Font.Name = "Aptos Narrow"

$human.Range("A2").Value = "ML"
$human.Range("A3").Value = "ML"

$human.PageSetup.Orientation = 1

$human.Range("A4").Select()
